# Insert a new weekly record at row 5 (above the existing row 5),
# shifting all subsequent data rows down by one. This mirrors the
# author's commit of adding a new "Fruta / hortaliza, semanal" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 5 (pushes rows 5..101 -> 6..102)
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new record's data
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44756
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 100114002
$ws.Cells.Item(5, 7).Value = "Camote"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 80
$ws.Cells.Item(5, 11).Value = 20000
$ws.Cells.Item(5, 12).Value = 20000
$ws.Cells.Item(5, 13).Value = 20000
$ws.Cells.Item(5, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 1000
$ws.Cells.Item(5, 17).Value = 20
$ws.Cells.Item(5, 18).Value = "Hortaliza"
